$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily refresh of the cryptos list (prices / 1h volume %) from GitHub Actions.
# Price cells (column D) are stored as text in this sheet; where the new
# price would otherwise be auto-recognized as a number by Excel, a leading
# apostrophe forces it to stay text (e.g. '''247.66' => the literal text
# "'247.66", which Excel stores as the text "247.66").
$ws.Range("D2").Value = '30.597.44'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '1.922.79'
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '''247.66'
$ws.Range("E5").Value = '  +2.99%  '
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").Value = '''0.4740'
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("D8").Value = '''0.2900'
$ws.Range("E8").Value = '  +0.90%  '
$ws.Range("D9").Value = '''0.06832'
$ws.Range("E9").Value = '  +3.84%  '
$ws.Range("D10").Value = '''105.55'
$ws.Range("E10").Value = '  -2.02%  '
$ws.Range("D11").Value = '''18.38'
$ws.Range("E11").Value = '  -3.62%  '
$ws.Range("D12").Value = '1.923.54'
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").Value = '''0.07711'
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("D14").Value = '''5.362'
$ws.Range("E14").Value = '  +4.08%  '
$ws.Range("D15").Value = '''0.6698'
$ws.Range("E15").Value = '  +1.54%  '
$ws.Range("D16").Value = '''291.36'
$ws.Range("E16").Value = '  -5.62%  '
$ws.Range("D17").Value = '30.611.44'
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("D18").Value = '''0.000007637'
$ws.Range("E18").Value = '  +1.48%  '
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("E20").Value = '  -0.53%  '
$ws.Range("D21").Value = '''5.547'
$ws.Range("E21").Value = '  +4.40%  '
$ws.Range("D22").Value = '2.174.17'
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("D24").Value = '''6.450'
$ws.Range("E24").Value = '  +2.47%  '
$ws.Range("D25").Value = '''9.527'
$ws.Range("E25").Value = '  +2.96%  '
$ws.Range("D26").Value = '''167.82'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = '''20.93'
$ws.Range("E27").Value = '  +3.33%  '
$ws.Range("D28").Value = '''2.127'
$ws.Range("E28").Value = '  +4.41%  '
$ws.Range("D29").Value = '''0.1075'
$ws.Range("E29").Value = '  -3.75%  '
$ws.Range("D30").Value = '''1.408'
$ws.Range("E30").Value = '  +3.63%  '
$ws.Range("D31").Value = '''4.181'
$ws.Range("E31").Value = '  +1.99%  '
$ws.Range("E32").Value = '  +3.10%  '
$ws.Range("D33").Value = '''0.05024'
$ws.Range("E33").Value = '  +0.28%  '
$ws.Range("E34").Value = '  -1.37%  '
$ws.Range("E35").Value = '  -0.50%  '
$ws.Range("D36").Value = '''0.02070'
$ws.Range("E36").Value = '  +6.14%  '
$ws.Range("D37").Value = '''0.9998'
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("D38").Value = '''2.725'
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("D39").Value = '''2.684'
$ws.Range("E39").Value = '  -0.58%  '
$ws.Range("D40").Value = '''111.47'
$ws.Range("E40").Value = '  +4.04%  '
$ws.Range("E41").Value = '  -0.79%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '''0.4436'
$ws.Range("E42").Value = '  +6.26%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '''0.8750'
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("D44").Value = '''5.898'
$ws.Range("E44").Value = '  +1.44%  '
$ws.Range("D45").Value = '''1.000'
$ws.Range("D46").Value = '''67.79'
$ws.Range("E46").Value = '  -3.81%  '
$ws.Range("D47").Value = '''7.307'
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("D48").Value = '''9.334'
$ws.Range("E48").Value = '  +0.31%  '
$ws.Range("D49").Value = '''48.21'
$ws.Range("E49").Value = '  +13.80%  '
$ws.Range("D50").Value = '''0.1246'
$ws.Range("E50").Value = '  +3.41%  '
$ws.Range("D51").Value = '''35.01'
$ws.Range("E51").Value = '  +0.31%  '
